# Sort the comma-separated "Recorded By" names in column G alphabetically
# using ordinal (ASCII, case-sensitive, uppercase-before-lowercase) order
# for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $list = New-Object System.Collections.Generic.List[string]
            foreach ($p in $parts) { [void]$list.Add($p) }
            [void]$list.Sort([System.StringComparer]::Ordinal)
            $newVal = [string]::Join(", ", $list)
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
